$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("hotel_info")

# Update hotel_info row 2 with review counts / rank data.
# Format as text first so the numeric-looking strings are stored as shared
# strings (t="s"), matching the source data's text-typed columns.
$dataRange = $ws1.Range("G2:I2")
$dataRange.NumberFormat = "@"

$ws1.Range("G2").Value = "5"
$ws1.Range("H2").Value = "462"
$ws1.Range("I2").Value = "5"
